{"js": "// Update the worksheet date header and all 25 \"NNN\u00d7N=\" practice problems\n// to the new day's values (2025-02-06 Thursday -> 2025-02-07 Friday, plus\n// new three-digit-by-one-digit multiplication problems).\nconst replacements = [\n  [\"2025-02-06 Thursday\", \"2025-02-07 Friday\"],\n  [\"564\u00d77=\", \"315\u00d75=\"],\n  [\"379\u00d73=\", \"218\u00d77=\"],\n  [\"298\u00d73=\", \"970\u00d72=\"],\n  [\"403\u00d76=\", \"767\u00d72=\"],\n  [\"481\u00d76=\", \"876\u00d79=\"],\n  [\"527\u00d72=\", \"724\u00d74=\"],\n  [\"740\u00d74=\", \"976\u00d73=\"],\n  [\"849\u00d78=\", \"323\u00d74=\"],\n  [\"426\u00d75=\", \"469\u00d72=\"],\n  [\"111\u00d76=\", \"745\u00d77=\"],\n  [\"723\u00d73=\", \"460\u00d78=\"],\n  [\"375\u00d72=\", \"468\u00d73=\"],\n  [\"966\u00d72=\", \"721\u00d74=\"],\n  [\"358\u00d79=\", \"317\u00d76=\"],\n  [\"575\u00d75=\", \"631\u00d77=\"],\n  [\"987\u00d72=\", \"956\u00d72=\"],\n  [\"175\u00d74=\", \"953\u00d76=\"],\n  [\"583\u00d76=\", \"161\u00d78=\"],\n  [\"488\u00d79=\", \"267\u00d77=\"],\n  [\"429\u00d76=\", \"498\u00d78=\"],\n  [\"901\u00d79=\", \"229\u00d74=\"],\n  [\"293\u00d75=\", \"744\u00d77=\"],\n  [\"707\u00d77=\", \"437\u00d74=\"],\n  [\"352\u00d76=\", \"435\u00d79=\"],\n  [\"710\u00d73=\", \"837\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and all 25 \"NNN\u00d7N=\" practice problems\n# to the new day's values (2025-02-06 Thursday -> 2025-02-07 Friday, plus\n# new three-digit-by-one-digit multiplication problems).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-06 Thursday\", \"2025-02-07 Friday\"),\n    @(\"564\u00d77=\", \"315\u00d75=\"),\n    @(\"379\u00d73=\", \"218\u00d77=\"),\n    @(\"298\u00d73=\", \"970\u00d72=\"),\n    @(\"403\u00d76=\", \"767\u00d72=\"),\n    @(\"481\u00d76=\", \"876\u00d79=\"),\n    @(\"527\u00d72=\", \"724\u00d74=\"),\n    @(\"740\u00d74=\", \"976\u00d73=\"),\n    @(\"849\u00d78=\", \"323\u00d74=\"),\n    @(\"426\u00d75=\", \"469\u00d72=\"),\n    @(\"111\u00d76=\", \"745\u00d77=\"),\n    @(\"723\u00d73=\", \"460\u00d78=\"),\n    @(\"375\u00d72=\", \"468\u00d73=\"),\n    @(\"966\u00d72=\", \"721\u00d74=\"),\n    @(\"358\u00d79=\", \"317\u00d76=\"),\n    @(\"575\u00d75=\", \"631\u00d77=\"),\n    @(\"987\u00d72=\", \"956\u00d72=\"),\n    @(\"175\u00d74=\", \"953\u00d76=\"),\n    @(\"583\u00d76=\", \"161\u00d78=\"),\n    @(\"488\u00d79=\", \"267\u00d77=\"),\n    @(\"429\u00d76=\", \"498\u00d78=\"),\n    @(\"901\u00d79=\", \"229\u00d74=\"),\n    @(\"293\u00d75=\", \"744\u00d77=\"),\n    @(\"707\u00d77=\", \"437\u00d74=\"),\n    @(\"352\u00d76=\", \"435\u00d79=\"),\n    @(\"710\u00d73=\", \"837\u00d79=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$true, [ref]$newText, [ref]2)\n}\n"}
